$wb = $excel.ActiveWorkbook

$updates = @{
    'Citywide Totals' = @{ 'K2'='7033'; 'J3'='8078'; 'K3'='7293'; 'E4'='2038'; 'F4'='1920'; 'H4'='1743'; 'K4'='1515'; 'K5'='521'; 'K6'='8018'; 'E7'='26044'; 'F7'='24113'; 'H7'='26056'; 'J7'='29307'; 'K7'='24380' }
    'By Neighborhood' = @{ 'K5'='66'; 'K7'='737'; 'K8'='1590'; 'K9'='112'; 'K10'='137'; 'K20'='597'; 'K21'='82'; 'K29'='1341'; 'K33'='1045'; 'K35'='37'; 'K36'='311'; 'K37'='824'; 'K41'='168'; 'K42'='898'; 'K47'='165'; 'K49'='135'; 'K51'='308'; 'K60'='140'; 'E63'='373'; 'F63'='204'; 'H63'='295'; 'K63'='68'; 'K65'='564'; 'K67'='950'; 'K68'='65'; 'K77'='164'; 'K78'='293'; 'K79'='595'; 'K80'='87'; 'K83'='521'; 'K85'='1119'; 'K88'='263'; 'K89'='362'; 'J90'='308'; 'K91'='289'; 'K94'='326'; 'K95'='397'; 'K97'='195'; 'K99'='412'; 'E101'='26044'; 'F101'='24113'; 'H101'='26056'; 'J101'='29307'; 'K101'='24380' }
    'Auburn Gresham' = @{ 'K3'='235'; 'K7'='737' }
    'Uptown' = @{ 'K6'='106'; 'K7'='362' }
    'South Shore' = @{ 'K2'='367'; 'K6'='275'; 'K7'='1119' }
    'Austin' = @{ 'K4'='90'; 'K6'='528'; 'K7'='1590' }
    'South Chicago' = @{ 'K3'='184'; 'K6'='121'; 'K7'='521' }
    'Garfield Park' = @{ 'K2'='261'; 'K3'='375'; 'K7'='1045' }
    'West Pullman' = @{ 'K4'='18'; 'K7'='397' }
    'Grand Crossing' = @{ 'K2'='238'; 'K7'='824' }
    'New City' = @{ 'K2'='185'; 'K4'='22'; 'K7'='564' }
    'Woodlawn' = @{ 'K3'='172'; 'K7'='412' }
    'North Lawndale' = @{ 'K3'='346'; 'K7'='950' }
    'Lincoln Park' = @{ 'K4'='13'; 'K6'='67'; 'K7'='135' }
    'Englewood' = @{ 'K2'='379'; 'K6'='391'; 'K7'='1341' }
    'Hermosa' = @{ 'K3'='36'; 'K7'='168' }
    'Humboldt Park' = @{ 'K3'='266'; 'K5'='15'; 'K7'='898' }
    'Avondale' = @{ 'K6'='61'; 'K7'='137' }
    'Rogers Park' = @{ 'K2'='87'; 'K6'='99'; 'K7'='293' }
    'Washington Park' = @{ 'K2'='74'; 'K3'='138'; 'K7'='289' }
    'Chinatown' = @{ 'K5'='4'; 'K7'='82' }
    'Roseland' = @{ 'K3'='191'; 'K6'='148'; 'K7'='595' }
    'Chicago Lawn' = @{ 'K4'='28'; 'K6'='164'; 'K7'='597' }
    'Grand Boulevard' = @{ 'K3'='97'; 'K5'='5'; 'K7'='311' }
    'West Loop' = @{ 'K6'='149'; 'K7'='326' }
    'Kenwood' = @{ 'K6'='51'; 'K7'='165' }
    'Gold Coast' = @{ 'K4'='5'; 'K7'='37' }
    'Avalon Park' = @{ 'K6'='29'; 'K7'='112' }
    'West Town' = @{ 'K3'='43'; 'K4'='7'; 'K6'='103'; 'K7'='195' }
    'United Center' = @{ 'K3'='81'; 'K7'='263' }
    'Armour Square' = @{ 'K6'='32'; 'K7'='66' }
    'Washington Heights' = @{ 'J3'='86'; 'J7'='308' }
    'Little Italy, UIC' = @{ 'K2'='83'; 'K7'='308' }
    'North Park' = @{ 'K6'='19'; 'K7'='65' }
    'Morgan Park' = @{ 'K2'='46'; 'K3'='43'; 'K7'='140' }
    'Riverdale' = @{ 'K3'='63'; 'K7'='164' }
    'Rush & Division' = @{ 'K4'='9'; 'K7'='87' }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = [double]$cellMap[$cellRef]
    }
}